$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 827 (shifts existing rows 827-863 down to 829-865)
$ws.Rows("827:828").Insert()

# New row 827: Piña, Caramelo, Primera, 12-unit box
$ws.Cells.Item(827,1).Value = 10
$ws.Cells.Item(827,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(827,3).Value = "La Araucanía"
$ws.Cells.Item(827,4).Value = 45267
$ws.Cells.Item(827,5).Value = 9
$ws.Cells.Item(827,6).Value = "Fruta"
$ws.Cells.Item(827,7).Value = 100108
$ws.Cells.Item(827,8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(827,9).Value = 100108005
$ws.Cells.Item(827,10).Value = "Piña"
$ws.Cells.Item(827,11).Value = "Caramelo"
$ws.Cells.Item(827,12).Value = "Primera"
$ws.Cells.Item(827,13).Value = 70
$ws.Cells.Item(827,14).Value = 25000
$ws.Cells.Item(827,15).Value = 26000
$ws.Cells.Item(827,16).Value = 25500
$ws.Cells.Item(827,17).Value = "`$/caja 12 unidades"
$ws.Cells.Item(827,18).Value = "Ecuador"
$ws.Cells.Item(827,19).Value = 2125
$ws.Cells.Item(827,20).Value = 12

# New row 828: Piña, Caramelo, Segunda, 14-unit box
$ws.Cells.Item(828,1).Value = 10
$ws.Cells.Item(828,2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(828,3).Value = "La Araucanía"
$ws.Cells.Item(828,4).Value = 45267
$ws.Cells.Item(828,5).Value = 9
$ws.Cells.Item(828,6).Value = "Fruta"
$ws.Cells.Item(828,7).Value = 100108
$ws.Cells.Item(828,8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(828,9).Value = 100108005
$ws.Cells.Item(828,10).Value = "Piña"
$ws.Cells.Item(828,11).Value = "Caramelo"
$ws.Cells.Item(828,12).Value = "Segunda"
$ws.Cells.Item(828,13).Value = 60
$ws.Cells.Item(828,14).Value = 26000
$ws.Cells.Item(828,15).Value = 26000
$ws.Cells.Item(828,16).Value = 26000
$ws.Cells.Item(828,17).Value = "`$/caja 14 unidades"
$ws.Cells.Item(828,18).Value = "Ecuador"
$ws.Cells.Item(828,19).Value = 1857
$ws.Cells.Item(828,20).Value = 14
